# Auto-generated script applying the Kraken_Profits market-data refresh
# to the ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G2=5489)
$ws.Cells.Item(2, 8).Value = 260.4  # H2: 260.6 -> 260.4
$ws.Cells.Item(2, 9).Value = 260.4  # I2: 315.75 -> 260.4
$ws.Cells.Item(2, 10).Value = 0  # J2: 40 -> 0
$ws.Cells.Item(2, 11).Value = 260.4  # K2: 315.75 -> 260.4
$ws.Cells.Item(2, 12).Value = 0  # L2: 40 -> 0
$ws.Cells.Item(2, 13).Value = -147.4  # M2: -202.75 -> -147.4
$ws.Cells.Item(2, 14).ClearContents()  # N2: -266 -> (removed)

# Row 19 (G19=7015)
$ws.Cells.Item(19, 8).Value = 1502.8  # H19: 1578.625 -> 1502.8
$ws.Cells.Item(19, 9).Value = 1374.5  # I19: 1499.6666 -> 1374.5
$ws.Cells.Item(19, 10).Value = 1588.3334  # J19: 1626 -> 1588.3334
$ws.Cells.Item(19, 11).Value = 1374.5  # K19: 1499.6666 -> 1374.5
$ws.Cells.Item(19, 12).Value = 1588.3334  # L19: 1626 -> 1588.3334
$ws.Cells.Item(19, 13).Value = -1199.5  # M19: -1324.6666 -> -1199.5
$ws.Cells.Item(19, 14).Value = -1938.3334  # N19: -1976 -> -1938.3334

# Row 28 (G28=27772)
$ws.Cells.Item(28, 8).Value = 693  # H28: 591.7778 -> 693
$ws.Cells.Item(28, 9).Value = 604.4286  # I28: 540 -> 604.4286
$ws.Cells.Item(28, 10).Value = 1003  # J28: 1006 -> 1003
$ws.Cells.Item(28, 11).Value = 604.4286  # K28: 540 -> 604.4286
$ws.Cells.Item(28, 12).Value = 1003  # L28: 1006 -> 1003
$ws.Cells.Item(28, 13).Value = -119.4286  # M28: -55 -> -119.4286
$ws.Cells.Item(28, 14).Value = -1973  # N28: -1976 -> -1973

# Row 29 (G29=4575)
$ws.Cells.Item(29, 8).Value = 785.2  # H29: 1330.9231 -> 785.2
$ws.Cells.Item(29, 9).Value = 88  # I29: 100.666664 -> 88
$ws.Cells.Item(29, 10).Value = 1250  # J29: 1700 -> 1250
$ws.Cells.Item(29, 11).Value = 264  # K29: 301.999992 -> 264
$ws.Cells.Item(29, 12).Value = 3750  # L29: 5100 -> 3750
$ws.Cells.Item(29, 13).Value = 17  # M29: -20.99999200000002 -> 17
$ws.Cells.Item(29, 14).Value = -4312  # N29: -5662 -> -4312

# Row 33 (G33=5512)
$ws.Cells.Item(33, 8).Value = 235  # H33: 203.14285 -> 235
$ws.Cells.Item(33, 9).Value = 275  # I33: 203.66667 -> 275
$ws.Cells.Item(33, 10).Value = 175  # J33: 200 -> 175
$ws.Cells.Item(33, 11).Value = 275  # K33: 203.66667 -> 275
$ws.Cells.Item(33, 12).Value = 175  # L33: 200 -> 175
$ws.Cells.Item(33, 13).Value = -46  # M33: 25.33332999999999 -> -46
$ws.Cells.Item(33, 14).Value = -633  # N33: -658 -> -633

# Row 43 (G43=5472)
$ws.Cells.Item(43, 8).Value = 2868.75  # H43: 2487.4167 -> 2868.75
$ws.Cells.Item(43, 9).Value = 2100  # I43: 1666.6666 -> 2100
$ws.Cells.Item(43, 10).Value = 3125  # J43: 2761 -> 3125
$ws.Cells.Item(43, 11).Value = 2100  # K43: 1666.6666 -> 2100
$ws.Cells.Item(43, 12).Value = 3125  # L43: 2761 -> 3125
$ws.Cells.Item(43, 13).Value = -2031  # M43: -1597.6666 -> -2031
$ws.Cells.Item(43, 14).Value = -3263  # N43: -2899 -> -3263

# Row 51 (G51=5486)
$ws.Cells.Item(51, 8).Value = 2000  # H51: 2500 -> 2000
$ws.Cells.Item(51, 9).Value = 2000  # I51: 2500 -> 2000
$ws.Cells.Item(51, 10).Value = 2000  # J51: 0 -> 2000
$ws.Cells.Item(51, 11).Value = 2000  # K51: 2500 -> 2000
$ws.Cells.Item(51, 12).Value = 2000  # L51: 0 -> 2000
$ws.Cells.Item(51, 13).Value = -1516  # M51: -2016 -> -1516
$ws.Cells.Item(51, 14).Value = -2968  # N51: None -> -2968

# Row 58 (G58=4606)
$ws.Cells.Item(58, 8).Value = 1466.5454  # H58: 1443.5714 -> 1466.5454
$ws.Cells.Item(58, 9).Value = 447.7143  # I58: 621 -> 447.7143
$ws.Cells.Item(58, 10).Value = 3249.5  # J58: 3500 -> 3249.5
$ws.Cells.Item(58, 11).Value = 1343.1429  # K58: 1863 -> 1343.1429
$ws.Cells.Item(58, 12).Value = 9748.5  # L58: 10500 -> 9748.5
$ws.Cells.Item(58, 13).Value = -1193.1429  # M58: -1713 -> -1193.1429
$ws.Cells.Item(58, 14).Value = -10048.5  # N58: -10800 -> -10048.5

# Row 88 (G88=12608)
$ws.Cells.Item(88, 8).Value = 3997.8333  # H88: 4164.5 -> 3997.8333
$ws.Cells.Item(88, 10).Value = 3332.6667  # J88: 3666 -> 3332.6667
$ws.Cells.Item(88, 12).Value = 3332.6667  # L88: 3666 -> 3332.6667
$ws.Cells.Item(88, 14).Value = -4144.6667  # N88: -4478 -> -4144.6667

# Row 91 (G91=12608)
$ws.Cells.Item(91, 8).Value = 3997.8333  # H91: 4164.5 -> 3997.8333
$ws.Cells.Item(91, 10).Value = 3332.6667  # J91: 3666 -> 3332.6667
$ws.Cells.Item(91, 12).Value = 3332.6667  # L91: 3666 -> 3332.6667
$ws.Cells.Item(91, 14).Value = -6140.6667  # N91: -6474 -> -6140.6667

# Row 98 (G98=36237)
$ws.Cells.Item(98, 8).Value = 1262.1818  # H98: 1157.4166 -> 1262.1818
$ws.Cells.Item(98, 9).Value = 1213.4  # I98: 1103.5454 -> 1213.4
$ws.Cells.Item(98, 11).Value = 1213.4  # K98: 1103.5454 -> 1213.4
$ws.Cells.Item(98, 13).Value = 284.5999999999999  # M98: 394.4546 -> 284.5999999999999

# Row 122 (G122=36237)
$ws.Cells.Item(122, 8).Value = 1262.1818  # H122: 1157.4166 -> 1262.1818
$ws.Cells.Item(122, 9).Value = 1213.4  # I122: 1103.5454 -> 1213.4
$ws.Cells.Item(122, 11).Value = 3640.2  # K122: 3310.6362 -> 3640.2
$ws.Cells.Item(122, 13).Value = -1190.2  # M122: -860.6361999999999 -> -1190.2

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (G94=19939)
$ws.Cells.Item(94, 8).Value = 6183.857  # H94: 5918 -> 6183.857
$ws.Cells.Item(94, 9).Value = 4996.75  # I94: 4993.6665 -> 4996.75
$ws.Cells.Item(94, 11).Value = 4996.75  # K94: 4993.6665 -> 4996.75
$ws.Cells.Item(94, 13).Value = -4545.75  # M94: -4542.6665 -> -4545.75

# Row 99 (G99=19943)
$ws.Cells.Item(99, 8).Value = 675  # H99: 400 -> 675
$ws.Cells.Item(99, 10).Value = 1500  # J99: 0 -> 1500
$ws.Cells.Item(99, 12).Value = 1500  # L99: 0 -> 1500
$ws.Cells.Item(99, 14).Value = -4496  # N99: None -> -4496

# Row 107 (G107=27706)
$ws.Cells.Item(107, 8).Value = 2110  # H107: 2110.5 -> 2110
$ws.Cells.Item(107, 9).Value = 2110  # I107: 2110.5 -> 2110
$ws.Cells.Item(107, 11).Value = 2110  # K107: 2110.5 -> 2110
$ws.Cells.Item(107, 13).Value = -190  # M107: -190.5 -> -190

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (G7=5361)
$ws.Cells.Item(7, 8).Value = 102.28571  # H7: 133.5 -> 102.28571
$ws.Cells.Item(7, 9).Value = 104.333336  # I7: 119.8 -> 104.333336
$ws.Cells.Item(7, 10).Value = 98.59999999999999  # J7: 156.33333 -> 98.59999999999999
$ws.Cells.Item(7, 11).Value = 104.333336  # K7: 119.8 -> 104.333336
$ws.Cells.Item(7, 12).Value = 98.59999999999999  # L7: 156.33333 -> 98.59999999999999
$ws.Cells.Item(7, 13).Value = 8.666663999999997  # M7: -6.799999999999997 -> 8.666663999999997
$ws.Cells.Item(7, 14).Value = -324.6  # N7: -382.33333 -> -324.6

# Row 62 (G62=12580)
$ws.Cells.Item(62, 8).Value = 800  # H62: 805 -> 800
$ws.Cells.Item(62, 9).Value = 800  # I62: 805 -> 800
$ws.Cells.Item(62, 10).Value = 800  # J62: 0 -> 800
$ws.Cells.Item(62, 11).Value = 800  # K62: 805 -> 800
$ws.Cells.Item(62, 12).Value = 800  # L62: 0 -> 800
$ws.Cells.Item(62, 13).Value = -176  # M62: -181 -> -176
$ws.Cells.Item(62, 14).Value = -2048  # N62: None -> -2048

# Row 65 (G65=12580)
$ws.Cells.Item(65, 8).Value = 800  # H65: 805 -> 800
$ws.Cells.Item(65, 9).Value = 800  # I65: 805 -> 800
$ws.Cells.Item(65, 10).Value = 800  # J65: 0 -> 800
$ws.Cells.Item(65, 11).Value = 4000  # K65: 4025 -> 4000
$ws.Cells.Item(65, 12).Value = 800  # L65: 0 -> 800
$ws.Cells.Item(65, 13).Value = -880  # M65: -905 -> -880
$ws.Cells.Item(65, 14).Value = -10240  # N65: None -> -10240

# Row 96 (G96=18193)
$ws.Cells.Item(96, 8).Value = 4703.2856  # H96: 3250 -> 4703.2856
$ws.Cells.Item(96, 10).Value = 4703.2856  # J96: 3250 -> 4703.2856
$ws.Cells.Item(96, 12).Value = 4703.2856  # L96: 3250 -> 4703.2856
$ws.Cells.Item(96, 14).Value = -10195.2856  # N96: -8742 -> -10195.2856

# Row 99 (G99=36198)
$ws.Cells.Item(99, 8).Value = 0  # H99: 2000 -> 0
$ws.Cells.Item(99, 9).Value = 0  # I99: 2000 -> 0
$ws.Cells.Item(99, 11).Value = 0  # K99: 2000 -> 0
$ws.Cells.Item(99, 13).ClearContents()  # M99: -502 -> (removed)

# Row 122 (G122=36196)
$ws.Cells.Item(122, 8).Value = 746  # H122: 617 -> 746
$ws.Cells.Item(122, 9).Value = 746  # I122: 617 -> 746
$ws.Cells.Item(122, 11).Value = 2238  # K122: 1851 -> 2238
$ws.Cells.Item(122, 13).Value = 212  # M122: 599 -> 212

# Row 126 (G126=36198)
$ws.Cells.Item(126, 8).Value = 0  # H126: 2000 -> 0
$ws.Cells.Item(126, 9).Value = 0  # I126: 2000 -> 0
$ws.Cells.Item(126, 11).Value = 0  # K126: 6000 -> 0
$ws.Cells.Item(126, 13).ClearContents()  # M126: -3530 -> (removed)

$ws = $wb.Worksheets.Item("CUL")
# Row 11 (G11=4745)
$ws.Cells.Item(11, 8).Value = 12500003  # H11: 50000000 -> 12500003
$ws.Cells.Item(11, 9).Value = 25000000  # I11: 50000000 -> 25000000
$ws.Cells.Item(11, 10).Value = 5  # J11: 0 -> 5
$ws.Cells.Item(11, 11).Value = 75000000  # K11: 150000000 -> 75000000
$ws.Cells.Item(11, 12).Value = 15  # L11: 0 -> 15
$ws.Cells.Item(11, 13).Value = -74999860  # M11: -149999860 -> -74999860
$ws.Cells.Item(11, 14).Value = -295  # N11: None -> -295

# Row 34 (G34=4749)
$ws.Cells.Item(34, 8).Value = 748.8570999999999  # H34: 829.875 -> 748.8570999999999
$ws.Cells.Item(34, 10).Value = 1161.75  # J34: 1208.8 -> 1161.75
$ws.Cells.Item(34, 12).Value = 3485.25  # L34: 3626.4 -> 3485.25
$ws.Cells.Item(34, 14).Value = -3653.25  # N34: -3794.4 -> -3653.25

# Row 39 (G39=4712)
$ws.Cells.Item(39, 8).Value = 2240.25  # H39: 1617.375 -> 2240.25
$ws.Cells.Item(39, 10).Value = 2000  # J39: 1329.6666 -> 2000
$ws.Cells.Item(39, 12).Value = 6000  # L39: 3988.9998 -> 6000
$ws.Cells.Item(39, 14).Value = -6588  # N39: -4576.9998 -> -6588

# Row 52 (G52=31902)
$ws.Cells.Item(52, 8).Value = 1000  # H52: 0 -> 1000
$ws.Cells.Item(52, 10).Value = 1000  # J52: 0 -> 1000
$ws.Cells.Item(52, 12).Value = 3000  # L52: 0 -> 3000
$ws.Cells.Item(52, 14).Value = -3532  # N52: None -> -3532

# Row 55 (G55=4733)
$ws.Cells.Item(55, 8).Value = 1293.2222  # H55: 1563.9 -> 1293.2222
$ws.Cells.Item(55, 10).Value = 0  # J55: 4000 -> 0
$ws.Cells.Item(55, 12).Value = 0  # L55: 12000 -> 0
$ws.Cells.Item(55, 14).ClearContents()  # N55: -12354 -> (removed)

# Row 86 (G86=12892)
$ws.Cells.Item(86, 8).Value = 213  # H86: 215 -> 213
$ws.Cells.Item(86, 9).Value = 213  # I86: 215 -> 213
$ws.Cells.Item(86, 11).Value = 639  # K86: 645 -> 639
$ws.Cells.Item(86, 13).Value = 547  # M86: 541 -> 547

# Row 89 (G89=12892)
$ws.Cells.Item(89, 8).Value = 213  # H89: 215 -> 213
$ws.Cells.Item(89, 9).Value = 213  # I89: 215 -> 213
$ws.Cells.Item(89, 11).Value = 1917  # K89: 1935 -> 1917
$ws.Cells.Item(89, 13).Value = 4011  # M89: 3993 -> 4011

$ws = $wb.Worksheets.Item("GSM")
# Row 46 (G46=2078)
$ws.Cells.Item(46, 8).Value = 0  # H46: 13285.714 -> 0
$ws.Cells.Item(46, 9).Value = 0  # I46: 4333.3335 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 20000 -> 0
$ws.Cells.Item(46, 11).Value = 0  # K46: 4333.3335 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 20000 -> 0
$ws.Cells.Item(46, 13).ClearContents()  # M46: -4177.3335 -> (removed)
$ws.Cells.Item(46, 14).ClearContents()  # N46: -20312 -> (removed)

# Row 133 (G133=41854)
$ws.Cells.Item(133, 8).Value = 99995  # H133: 0 -> 99995
$ws.Cells.Item(133, 10).Value = 99995  # J133: 0 -> 99995
$ws.Cells.Item(133, 12).Value = 99995  # L133: 0 -> 99995
$ws.Cells.Item(133, 14).Value = -110115  # N133: None -> -110115

# Row 134 (G134=42064)
$ws.Cells.Item(134, 8).Value = 99000  # H134: 99666.664 -> 99000
$ws.Cells.Item(134, 10).Value = 99000  # J134: 99666.664 -> 99000
$ws.Cells.Item(134, 12).Value = 297000  # L134: 298999.992 -> 297000
$ws.Cells.Item(134, 14).Value = -302070  # N134: -304069.992 -> -302070

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (G16=5289)
$ws.Cells.Item(16, 8).Value = 3111  # H16: 4000 -> 3111
$ws.Cells.Item(16, 9).Value = 3111  # I16: 4000 -> 3111
$ws.Cells.Item(16, 11).Value = 3111  # K16: 4000 -> 3111
$ws.Cells.Item(16, 13).Value = -2941  # M16: -3830 -> -2941

# Row 22 (G22=5277)
$ws.Cells.Item(22, 8).Value = 2133.3333  # H22: 2254.2856 -> 2133.3333
$ws.Cells.Item(22, 9).Value = 1400  # I22: 2195 -> 1400
$ws.Cells.Item(22, 10).Value = 2500  # J22: 2278 -> 2500
$ws.Cells.Item(22, 11).Value = 1400  # K22: 2195 -> 1400
$ws.Cells.Item(22, 12).Value = 2500  # L22: 2278 -> 2500
$ws.Cells.Item(22, 13).Value = -1105  # M22: -1900 -> -1105
$ws.Cells.Item(22, 14).Value = -3090  # N22: -2868 -> -3090

# Row 27 (G27=5277)
$ws.Cells.Item(27, 8).Value = 2133.3333  # H27: 2254.2856 -> 2133.3333
$ws.Cells.Item(27, 9).Value = 1400  # I27: 2195 -> 1400
$ws.Cells.Item(27, 10).Value = 2500  # J27: 2278 -> 2500
$ws.Cells.Item(27, 11).Value = 1400  # K27: 2195 -> 1400
$ws.Cells.Item(27, 12).Value = 2500  # L27: 2278 -> 2500
$ws.Cells.Item(27, 13).Value = -1293  # M27: -2088 -> -1293
$ws.Cells.Item(27, 14).Value = -2714  # N27: -2492 -> -2714

# Row 40 (G40=36248)
$ws.Cells.Item(40, 8).Value = 5548  # H40: 6638 -> 5548
$ws.Cells.Item(40, 9).Value = 5548  # I40: 6638 -> 5548
$ws.Cells.Item(40, 11).Value = 5548  # K40: 6638 -> 5548
$ws.Cells.Item(40, 13).Value = -5412  # M40: -6502 -> -5412

# Row 68 (G68=12563)
$ws.Cells.Item(68, 8).Value = 2187.5  # H68: 2214.2856 -> 2187.5
$ws.Cells.Item(68, 10).Value = 2750  # J68: 2375 -> 2750
$ws.Cells.Item(68, 12).Value = 2750  # L68: 2375 -> 2750
$ws.Cells.Item(68, 14).Value = -4248  # N68: -3873 -> -4248

# Row 71 (G71=12563)
$ws.Cells.Item(71, 8).Value = 2187.5  # H71: 2214.2856 -> 2187.5
$ws.Cells.Item(71, 10).Value = 2750  # J71: 2375 -> 2750
$ws.Cells.Item(71, 12).Value = 13750  # L71: 11875 -> 13750
$ws.Cells.Item(71, 14).Value = -21238  # N71: -19363 -> -21238

# Row 103 (G103=18526)
$ws.Cells.Item(103, 8).Value = 30000  # H103: 29999.5 -> 30000
$ws.Cells.Item(103, 10).Value = 30000  # J103: 29999.5 -> 30000
$ws.Cells.Item(103, 12).Value = 30000  # L103: 29999.5 -> 30000
$ws.Cells.Item(103, 14).Value = -32344  # N103: -32343.5 -> -32344

# Row 122 (G122=36247)
$ws.Cells.Item(122, 8).Value = 4394.5  # H122: 4394.8335 -> 4394.5
$ws.Cells.Item(122, 9).Value = 4473.4  # I122: 4473.8 -> 4473.4
$ws.Cells.Item(122, 11).Value = 13420.2  # K122: 13421.4 -> 13420.2
$ws.Cells.Item(122, 13).Value = -10970.2  # M122: -10971.4 -> -10970.2

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (G62=12589)
$ws.Cells.Item(62, 8).Value = 2618.8  # H62: 2898.5 -> 2618.8
$ws.Cells.Item(62, 9).Value = 2523.5  # I62: 2864.6667 -> 2523.5
$ws.Cells.Item(62, 11).Value = 2523.5  # K62: 2864.6667 -> 2523.5
$ws.Cells.Item(62, 13).Value = -1899.5  # M62: -2240.6667 -> -1899.5

# Row 65 (G65=12589)
$ws.Cells.Item(65, 8).Value = 2618.8  # H65: 2898.5 -> 2618.8
$ws.Cells.Item(65, 9).Value = 2523.5  # I65: 2864.6667 -> 2523.5
$ws.Cells.Item(65, 11).Value = 12617.5  # K65: 14323.3335 -> 12617.5
$ws.Cells.Item(65, 13).Value = -9497.5  # M65: -11203.3335 -> -9497.5

# Row 122 (G122=36208)
$ws.Cells.Item(122, 8).Value = 2016.1875  # H122: 2143 -> 2016.1875
$ws.Cells.Item(122, 9).Value = 2003.8572  # I122: 2149.2307 -> 2003.8572
$ws.Cells.Item(122, 11).Value = 6011.571599999999  # K122: 6447.6921 -> 6011.571599999999
$ws.Cells.Item(122, 13).Value = -3561.571599999999  # M122: -3997.6921 -> -3561.571599999999
